$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before column D (shifts D:K -> E:L)
$ws.Range("D1").EntireColumn.Insert(1, 1)

$ws.Range("D7").Value = 43465

Write-Output "done"
